# "Added 1 more test case espn"
# Rename the "OverStock" sheet to "YearResult" and populate it with a
# column of season labels (a new test-data sheet), then make it the
# active/selected sheet - matching the committed workbook state.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OverStock")
$ws.Name = "YearResult"

$values = @(
    "YearResult",
    "2021-22",
    "2020-21",
    "2019-20",
    "2018-19",
    "2017-18",
    "2016-17",
    "2015-16",
    "2014-15",
    "2013-14",
    "2012-13",
    "2011-12",
    "2010-11",
    "2009-10",
    "2008-09",
    "2007-08",
    "2006-07",
    "2005-06",
    "2004-05",
    "2003-04",
    "2002-03"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# This sheet becomes the selected/active tab in the saved workbook.
$ws.Activate()
